$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 3: "The quality control ... examined manually. " -> "... examined manually, i.e. extra step is needed."
$para3 = $tr.Paragraphs(3,1)
$run3 = $para3.Runs(1,1)
$run3.Text = "The quality control can be controlled totally on visual classification. In which it can be either binary classifier, so the classifier can detect the existence/absence of the defect. But the defect will be then examined manually, i.e. extra step is needed."

# Paragraph 5: "Business rule ... stop line/machine. " -> append " Or he will simply git rid of the defective device. "
$para5 = $tr.Paragraphs(5,1)
$run5 = $para5.Runs(1,1)
$run5.Text = "Business rule are then be implied to see what you will do with the defective devices. After discovery, what the quality engineer will do ? Will he count the number of defectives and if increased on threshold he will stop line/machine? Or he will calculate some probability and if the probability of defective > threshold he will stop line/machine. Or he will simply git rid of the defective device. "

# Add a brand-new paragraph at the end of the text frame.
$null = $tr.InsertAfter("`rOne concern for multiple classes is the existence of new class of defectives. We shall train the model again")
